$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29, shifting the existing rows 29-50 down to 30-51.
$ws.Rows.Item(29).Insert()

# Populate the new row 29 with this week's record (same market/category metadata,
# new date + volume + price figures).
$ws.Cells.Item(29, 1).Value = 10
$ws.Cells.Item(29, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(29, 3).Value = "La Araucanía"
$ws.Cells.Item(29, 4).Value = 44762
$ws.Cells.Item(29, 5).Value = 9
$ws.Cells.Item(29, 6).Value = 100112010
$ws.Cells.Item(29, 7).Value = "Achicoria"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 85
$ws.Cells.Item(29, 11).Value = 11000
$ws.Cells.Item(29, 12).Value = 11000
$ws.Cells.Item(29, 13).Value = 11000
$ws.Cells.Item(29, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(29, 15).Value = "Región Metropolitana"
$ws.Cells.Item(29, 16).Value = 611
$ws.Cells.Item(29, 17).Value = 18
$ws.Cells.Item(29, 18).Value = "Hortaliza"
